$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '24.973.30', '  -3.84%  ')
    ,@(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.639.19', '  -6.11%  ')
    ,@(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '0.9965', '  -0.34%  ')
    ,@(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '233.82', '  -6.03%  ')
    ,@(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.9989', '  -0.10%  ')
    ,@(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4757', '  -6.23%  ')
    ,@(8, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '39.24', '  -3.72%  ')
    ,@(9, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.2586', '  -6.13%  ')
    ,@(10, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.06091', '  -1.75%  ')
    ,@(11, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07014', '  -3.24%  ')
    ,@(12, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.640.90', '  -6.04%  ')
    ,@(13, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '14.68', '  -3.36%  ')
    ,@(14, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.6008', '  -8.17%  ')
    ,@(15, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '4.344', '  -7.05%  ')
    ,@(16, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '73.59', '  -5.41%  ')
    ,@(17, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.9990', '  -0.11%  ')
    ,@(18, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '0.9972', '  -0.24%  ')
    ,@(19, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '24.961.51', '  -3.92%  ')
    ,@(20, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000006586', '  -3.97%  ')
    ,@(21, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '11.20', '  -5.72%  ')
    ,@(22, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '1.851.26', '  -5.96%  ')
    ,@(23, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '4.371', '  -1.33%  ')
    ,@(24, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '8.572', '  -2.03%  ')
    ,@(25, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '5.267', '  -2.34%  ')
    ,@(26, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '133.47', '  -2.37%  ')
    ,@(27, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '14.84', '  -2.83%  ')
    ,@(28, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.383', '  -8.84%  ')
    ,@(29, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '103.54', '  -2.11%  ')
    ,@(30, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.639', '  -8.19%  ')
    ,@(31, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '3.959', '  +2.35%  ')
    ,@(32, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.07712', '  -5.97%  ')
    ,@(33, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '3.550', '  -2.88%  ')
    ,@(34, 'Frax', 'https://coinranking.com/coin/KfWtaeV1W+frax-frax', '0.9981', '  -0.08%  ')
    ,@(35, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.04303', '  -7.96%  ')
    ,@(36, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.589', '  -2.51%  ')
    ,@(37, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '0.9266', '  -7.23%  ')
    ,@(38, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.5829', '  -5.58%  ')
    ,@(39, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.538', '  -7.82%  ')
    ,@(40, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01531', '  -5.21%  ')
    ,@(41, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '0.9979', '  -0.18%  ')
    ,@(42, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.8167', '  +6.59%  ')
    ,@(43, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '97.86', '  -2.91%  ')
    ,@(44, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '1.769', '  -8.50%  ')
    ,@(45, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.3698', '  -5.84%  ')
    ,@(46, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '4.705', '  -6.05%  ')
    ,@(47, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1088', '  -5.72%  ')
    ,@(48, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '6.066', '  -4.47%  ')
    ,@(49, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.05193', '  -2.78%  ')
    ,@(50, 'Elrond', 'https://coinranking.com/coin/omwkOTglq+elrond-egld', '29.46', '  -4.13%  ')
    ,@(51, 'TrueUSD', 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd', '0.9985', '  -0.34%  ')
)

foreach ($r in $data) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 2).Value = $r[1]
    $ws.Cells.Item($rowNum, 3).Value = $r[2]
    $ws.Cells.Item($rowNum, 4).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 4).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Style = "Normal"
    $ws.Cells.Item($rowNum, 5).Value = $r[4]
}

Write-Output "done"